$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 15572.333
$ws.Range("J21").Value = 23000
$ws.Range("L21").Value = 23000
$ws.Range("N21").Value = -23936

$ws.Range("H23").Value = 15572.333
$ws.Range("J23").Value = 23000
$ws.Range("L23").Value = 23000
$ws.Range("N23").Value = -23468

$ws.Range("H106").Value = 35000
$ws.Range("I106").Value = 35000
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 35000
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = -34369
$ws.Range("N106").ClearContents()

$ws.Range("H112").Value = 1621.2188
$ws.Range("J112").Value = 1657.3871
$ws.Range("L112").Value = 4972.1613
$ws.Range("N112").Value = -7188.1613

$ws.Range("H132").Value = 4493.871
$ws.Range("I132").Value = 3986.818
$ws.Range("J132").Value = 5733.3335
$ws.Range("K132").Value = 11960.454
$ws.Range("L132").Value = 17200.0005
$ws.Range("M132").Value = -9430.454000000002
$ws.Range("N132").Value = -22260.0005

$ws.Range("H135").Value = 18520148
$ws.Range("I135").Value = 1693.1923
$ws.Range("K135").Value = 15238.7307
$ws.Range("M135").Value = -12703.7307

$ws.Range("H137").Value = 18182788
$ws.Range("I137").Value = 851.44446
$ws.Range("K137").Value = 2554.33338
$ws.Range("M137").Value = -4.333380000000034

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11801.976
$ws.Range("I32").Value = 11397.017
$ws.Range("J32").Value = 12858.392
$ws.Range("K32").Value = 11397.017
$ws.Range("L32").Value = 12858.392
$ws.Range("M32").Value = -11110.017
$ws.Range("N32").Value = -13432.392

$ws.Range("H45").Value = 1196
$ws.Range("I45").Value = 1105.4117
$ws.Range("J45").Value = 1350
$ws.Range("K45").Value = 1105.4117
$ws.Range("L45").Value = 1350
$ws.Range("M45").Value = -728.4117000000001
$ws.Range("N45").Value = -2104

$ws.Range("H61").Value = 15153697
$ws.Range("I61").Value = 16131145
$ws.Range("K61").Value = 16131145
$ws.Range("M61").Value = -16130933

$ws.Range("H122").Value = 24289.143
$ws.Range("I122").Value = 28004.8
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 84014.39999999999
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -81564.39999999999
$ws.Range("N122").Value = -49900

$ws.Range("H133").Value = 48190.25
$ws.Range("J133").Value = 48190.25
$ws.Range("L133").Value = 48190.25
$ws.Range("N133").Value = -53250.25

$ws.Range("H136").Value = 15153697
$ws.Range("I136").Value = 16131145
$ws.Range("K136").Value = 48393435
$ws.Range("M136").Value = -48390885

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 954.5454999999999
$ws.Range("I94").Value = 867.3684
$ws.Range("K94").Value = 867.3684
$ws.Range("M94").Value = -416.3684

$ws.Range("H105").Value = 3219.0454
$ws.Range("I105").Value = 1718
$ws.Range("K105").Value = 1718
$ws.Range("M105").Value = 29

$ws.Range("H134").Value = 4596.2144
$ws.Range("I134").Value = 2797.1667
$ws.Range("J134").Value = 5945.5
$ws.Range("K134").Value = 8391.500100000001
$ws.Range("L134").Value = 17836.5
$ws.Range("M134").Value = -5856.500100000001
$ws.Range("N134").Value = -22906.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 45999.5
$ws.Range("J20").Value = 45999.5
$ws.Range("L20").Value = 45999.5
$ws.Range("N20").Value = -46471.5

$ws.Range("H30").Value = 45999.5
$ws.Range("J30").Value = 45999.5
$ws.Range("L30").Value = 45999.5
$ws.Range("N30").Value = -46181.5

$ws.Range("H74").Value = 13599.4
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 13599.4
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 13599.4
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -15347.4

$ws.Range("H77").Value = 13599.4
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 13599.4
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 40798.2
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -49534.2

$ws.Range("H128").Value = 45999.5
$ws.Range("J128").Value = 45999.5
$ws.Range("L128").Value = 45999.5
$ws.Range("N128").Value = -55959.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 29.947369
$ws.Range("I12").Value = 35
$ws.Range("J12").Value = 27
$ws.Range("K12").Value = 105
$ws.Range("L12").Value = 81
$ws.Range("M12").Value = 68
$ws.Range("N12").Value = -427

$ws.Range("H104").Value = 6865
$ws.Range("I104").Value = 2026
$ws.Range("J104").Value = 8074.75
$ws.Range("K104").Value = 6078
$ws.Range("L104").Value = 24224.25
$ws.Range("M104").Value = -3457
$ws.Range("N104").Value = -29466.25

$ws.Range("H113").Value = 819.625
$ws.Range("I113").Value = 596.44446
$ws.Range("J113").Value = 953.5333000000001
$ws.Range("K113").Value = 1789.33338
$ws.Range("L113").Value = 2860.5999
$ws.Range("M113").Value = 380.66662
$ws.Range("N113").Value = -7200.5999

$ws.Range("H118").Value = 1223.9166
$ws.Range("J118").Value = 1447.6666
$ws.Range("L118").Value = 4342.9998
$ws.Range("N118").Value = -6828.9998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 28625
$ws.Range("I70").Value = 69000
$ws.Range("J70").Value = 4400
$ws.Range("K70").Value = 69000
$ws.Range("L70").Value = 4400
$ws.Range("M70").Value = -68730
$ws.Range("N70").Value = -4940

$ws.Range("H73").Value = 28625
$ws.Range("I73").Value = 69000
$ws.Range("J73").Value = 4400
$ws.Range("K73").Value = 69000
$ws.Range("L73").Value = 4400
$ws.Range("M73").Value = -68064
$ws.Range("N73").Value = -6272

$ws.Range("H132").Value = 5495.488
$ws.Range("I132").Value = 5539.3335
$ws.Range("J132").Value = 5179.8
$ws.Range("K132").Value = 16618.0005
$ws.Range("L132").Value = 15539.4
$ws.Range("M132").Value = -14088.0005
$ws.Range("N132").Value = -20599.4

$ws.Range("H140").Value = 45780
$ws.Range("J140").Value = 45780
$ws.Range("L140").Value = 45780
$ws.Range("N140").Value = -56140

$ws.Range("H141").Value = 399114.5
$ws.Range("J141").Value = 399114.5
$ws.Range("L141").Value = 399114.5
$ws.Range("N141").Value = -409474.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1672.1538
$ws.Range("I61").Value = 1354
$ws.Range("J61").Value = 2732.6667
$ws.Range("K61").Value = 1354
$ws.Range("L61").Value = 2732.6667
$ws.Range("M61").Value = -1152
$ws.Range("N61").Value = -3136.6667

$ws.Range("H113").Value = 1672.1538
$ws.Range("I113").Value = 1354
$ws.Range("J113").Value = 2732.6667
$ws.Range("K113").Value = 1354
$ws.Range("L113").Value = 2732.6667
$ws.Range("M113").Value = 816
$ws.Range("N113").Value = -7072.6667

$ws.Range("H132").Value = 19236528
$ws.Range("I132").Value = 6138.8237
$ws.Range("J132").Value = 55560600
$ws.Range("K132").Value = 18416.4711
$ws.Range("L132").Value = 166681800
$ws.Range("M132").Value = -15886.4711
$ws.Range("N132").Value = -166686860

$ws.Range("H133").Value = 55000
$ws.Range("J133").Value = 55000
$ws.Range("L133").Value = 55000
$ws.Range("N133").Value = -60060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 783.6070999999999
$ws.Range("I136").Value = 783.6070999999999
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2350.8213
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 199.1787000000004
$ws.Range("N136").ClearContents()
